$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(11, 8).Value = 484.86667  # ALC!H11 (was 560)
$ws.Cells.Item(11, 9).Value = 484.86667  # ALC!I11 (was 560)
$ws.Cells.Item(11, 11).Value = 484.86667  # ALC!K11 (was 560)
$ws.Cells.Item(11, 13).Value = -344.86667  # ALC!M11 (was -420)

$ws.Cells.Item(40, 8).Value = 4500  # ALC!H40 (was 2475)
$ws.Cells.Item(40, 9).Value = 2000  # ALC!I40 (was 2475)
$ws.Cells.Item(40, 10).Value = 7000  # ALC!J40 (was 0)
$ws.Cells.Item(40, 11).Value = 2000  # ALC!K40 (was 2475)
$ws.Cells.Item(40, 12).Value = 7000  # ALC!L40 (was 0)
$ws.Cells.Item(40, 13).Value = -1825  # ALC!M40 (was -2300)
$ws.Cells.Item(40, 14).Value = -7350  # ALC!N40 (was ADD)

$ws.Cells.Item(132, 8).Value = 1717.75  # ALC!H132 (was 1722.3125)
$ws.Cells.Item(132, 9).Value = 1781.3572  # ALC!I132 (was 1786.5714)
$ws.Cells.Item(132, 11).Value = 5344.071599999999  # ALC!K132 (was 5359.7142)
$ws.Cells.Item(132, 13).Value = -2814.071599999999  # ALC!M132 (was -2829.7142)

$ws.Cells.Item(141, 8).Value = 2631.3635  # ALC!H141 (was 2813.182)
$ws.Cells.Item(141, 9).Value = 2631.3635  # ALC!I141 (was 2813.182)
$ws.Cells.Item(141, 11).Value = 7894.0905  # ALC!K141 (was 8439.545999999998)
$ws.Cells.Item(141, 13).Value = -2714.0905  # ALC!M141 (was -3259.545999999998)

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(112, 8).Value = 16372  # ARM!H112 (was 16183.5)
$ws.Cells.Item(112, 10).Value = 16372  # ARM!J112 (was 16183.5)
$ws.Cells.Item(112, 12).Value = 16372  # ARM!L112 (was 16183.5)
$ws.Cells.Item(112, 14).Value = -19326  # ARM!N112 (was -19137.5)

$ws.Cells.Item(122, 8).Value = 0  # ARM!H122 (was 3005.75)
$ws.Cells.Item(122, 9).Value = 0  # ARM!I122 (was 3005.75)
$ws.Cells.Item(122, 11).Value = 0  # ARM!K122 (was 9017.25)
$ws.Cells.Item(122, 13).ClearContents()  # ARM!M122

$ws.Cells.Item(132, 8).Value = 2861.4285  # ARM!H132 (was 3019.4)
$ws.Cells.Item(132, 9).Value = 2767.8948  # ARM!I132 (was 2938.2222)
$ws.Cells.Item(132, 11).Value = 8303.6844  # ARM!K132 (was 8814.6666)
$ws.Cells.Item(132, 13).Value = -5773.6844  # ARM!M132 (was -6284.6666)

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 1197.7812  # BSM!H134 (was 1230.8438)
$ws.Cells.Item(134, 9).Value = 1197.7812  # BSM!I134 (was 1230.8438)
$ws.Cells.Item(134, 11).Value = 3593.3436  # BSM!K134 (was 3692.5314)
$ws.Cells.Item(134, 13).Value = -1058.3436  # BSM!M134 (was -1157.5314)

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1708.9584  # CRP!H16 (was 1721.0834)
$ws.Cells.Item(16, 9).Value = 1841.7142  # CRP!I16 (was 1771.6666)
$ws.Cells.Item(16, 10).Value = 1523.1  # CRP!J16 (was 1636.7778)
$ws.Cells.Item(16, 11).Value = 1841.7142  # CRP!K16 (was 1771.6666)
$ws.Cells.Item(16, 12).Value = 1523.1  # CRP!L16 (was 1636.7778)
$ws.Cells.Item(16, 13).Value = -1554.7142  # CRP!M16 (was -1484.6666)
$ws.Cells.Item(16, 14).Value = -2097.1  # CRP!N16 (was -2210.7778)

$ws.Cells.Item(22, 8).Value = 581.44446  # CRP!H22 (was 668.1429000000001)
$ws.Cells.Item(22, 9).Value = 404.125  # CRP!I22 (was 446.16666)
$ws.Cells.Item(22, 11).Value = 404.125  # CRP!K22 (was 446.16666)
$ws.Cells.Item(22, 13).Value = -54.125  # CRP!M22 (was -96.16665999999998)

$ws.Cells.Item(107, 8).Value = 1615.2174  # CRP!H107 (was 1885.5)
$ws.Cells.Item(107, 9).Value = 1478.6154  # CRP!I107 (was 1841.6923)
$ws.Cells.Item(107, 10).Value = 1792.8  # CRP!J107 (was 1948.7778)
$ws.Cells.Item(107, 11).Value = 1478.6154  # CRP!K107 (was 1841.6923)
$ws.Cells.Item(107, 12).Value = 1792.8  # CRP!L107 (was 1948.7778)
$ws.Cells.Item(107, 13).Value = 441.3846000000001  # CRP!M107 (was 78.30770000000007)
$ws.Cells.Item(107, 14).Value = -5632.8  # CRP!N107 (was -5788.7778)

$ws.Cells.Item(113, 8).Value = 1708.9584  # CRP!H113 (was 1721.0834)
$ws.Cells.Item(113, 9).Value = 1841.7142  # CRP!I113 (was 1771.6666)
$ws.Cells.Item(113, 10).Value = 1523.1  # CRP!J113 (was 1636.7778)
$ws.Cells.Item(113, 11).Value = 1841.7142  # CRP!K113 (was 1771.6666)
$ws.Cells.Item(113, 12).Value = 1523.1  # CRP!L113 (was 1636.7778)
$ws.Cells.Item(113, 13).Value = 328.2858000000001  # CRP!M113 (was 398.3334)
$ws.Cells.Item(113, 14).Value = -5863.1  # CRP!N113 (was -5976.7778)

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 1734.9412  # CUL!H34 (was 1708.5294)
$ws.Cells.Item(34, 9).Value = 219.83333  # CUL!I34 (was 202.57143)
$ws.Cells.Item(34, 10).Value = 2561.3635  # CUL!J34 (was 2762.7)
$ws.Cells.Item(34, 11).Value = 659.49999  # CUL!K34 (was 607.71429)
$ws.Cells.Item(34, 12).Value = 7684.0905  # CUL!L34 (was 8288.099999999999)
$ws.Cells.Item(34, 13).Value = -575.49999  # CUL!M34 (was -523.71429)
$ws.Cells.Item(34, 14).Value = -7852.0905  # CUL!N34 (was -8456.099999999999)

$ws.Cells.Item(51, 8).Value = 4999  # CUL!H51 (was 2599.5)
$ws.Cells.Item(51, 9).Value = 4999  # CUL!I51 (was 2599.5)
$ws.Cells.Item(51, 11).Value = 14997  # CUL!K51 (was 7798.5)
$ws.Cells.Item(51, 13).Value = -14537  # CUL!M51 (was -7338.5)

$ws.Cells.Item(131, 8).Value = 2126  # CUL!H131 (was 2055)

$ws.Cells.Item(140, 8).Value = 2454.9  # CUL!H140 (was 2331.6365)
$ws.Cells.Item(140, 9).Value = 2454.9  # CUL!I140 (was 2331.6365)
$ws.Cells.Item(140, 11).Value = 7364.700000000001  # CUL!K140 (was 6994.9095)
$ws.Cells.Item(140, 13).Value = -2184.700000000001  # CUL!M140 (was -1814.9095)

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 583.5625  # GSM!H2 (was 505.5)
$ws.Cells.Item(2, 9).Value = 246.08333  # GSM!I2 (was 218.92857)
$ws.Cells.Item(2, 10).Value = 1596  # GSM!J2 (was 1007)
$ws.Cells.Item(2, 11).Value = 246.08333  # GSM!K2 (was 218.92857)
$ws.Cells.Item(2, 12).Value = 1596  # GSM!L2 (was 1007)
$ws.Cells.Item(2, 13).Value = -133.08333  # GSM!M2 (was -105.92857)
$ws.Cells.Item(2, 14).Value = -1822  # GSM!N2 (was -1233)

$ws.Cells.Item(70, 8).Value = 9499.166999999999  # GSM!H70 (was 0)
$ws.Cells.Item(70, 9).Value = 9499.166999999999  # GSM!I70 (was 0)
$ws.Cells.Item(70, 11).Value = 9499.166999999999  # GSM!K70 (was 0)
$ws.Cells.Item(70, 13).Value = -9229.166999999999  # GSM!M70 (was ADD)

$ws.Cells.Item(73, 8).Value = 9499.166999999999  # GSM!H73 (was 0)
$ws.Cells.Item(73, 9).Value = 9499.166999999999  # GSM!I73 (was 0)
$ws.Cells.Item(73, 11).Value = 9499.166999999999  # GSM!K73 (was 0)
$ws.Cells.Item(73, 13).Value = -8563.166999999999  # GSM!M73 (was ADD)

$ws.Cells.Item(86, 8).Value = 25000  # GSM!H86 (was 24714.334)
$ws.Cells.Item(86, 10).Value = 25000  # GSM!J86 (was 24714.334)
$ws.Cells.Item(86, 12).Value = 25000  # GSM!L86 (was 24714.334)
$ws.Cells.Item(86, 14).Value = -27372  # GSM!N86 (was -27086.334)

$ws.Cells.Item(89, 8).Value = 25000  # GSM!H89 (was 24714.334)
$ws.Cells.Item(89, 10).Value = 25000  # GSM!J89 (was 24714.334)
$ws.Cells.Item(89, 12).Value = 75000  # GSM!L89 (was 74143.00199999999)
$ws.Cells.Item(89, 14).Value = -86856  # GSM!N89 (was -85999.00199999999)

$ws.Cells.Item(97, 8).Value = 27513.2  # GSM!H97 (was 28640.834)
$ws.Cells.Item(97, 9).Value = 47441.43  # GSM!I97 (was 51056.152)
$ws.Cells.Item(97, 11).Value = 47441.43  # GSM!K97 (was 51056.152)
$ws.Cells.Item(97, 13).Value = -46945.43  # GSM!M97 (was -50560.152)

$ws.Cells.Item(122, 8).Value = 2082.9092  # GSM!H122 (was 1992.25)
$ws.Cells.Item(122, 9).Value = 2082.9092  # GSM!I122 (was 1992.25)
$ws.Cells.Item(122, 11).Value = 6248.7276  # GSM!K122 (was 5976.75)
$ws.Cells.Item(122, 13).Value = -3798.7276  # GSM!M122 (was -3526.75)

$ws.Cells.Item(124, 8).Value = 151080  # GSM!H124 (was 0)
$ws.Cells.Item(124, 10).Value = 151080  # GSM!J124 (was 0)
$ws.Cells.Item(124, 12).Value = 151080  # GSM!L124 (was 0)
$ws.Cells.Item(124, 14).Value = -160900  # GSM!N124 (was ADD)

$ws.Cells.Item(126, 8).Value = 1532.0588  # GSM!H126 (was 1404.75)
$ws.Cells.Item(126, 9).Value = 1202  # GSM!I126 (was 1072.3334)
$ws.Cells.Item(126, 11).Value = 3606  # GSM!K126 (was 3217.0002)
$ws.Cells.Item(126, 13).Value = -1136  # GSM!M126 (was -747.0001999999999)

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(29, 8).Value = 0  # LTW!H29 (was 20000)
$ws.Cells.Item(29, 10).Value = 0  # LTW!J29 (was 20000)
$ws.Cells.Item(29, 12).Value = 0  # LTW!L29 (was 20000)
$ws.Cells.Item(29, 14).ClearContents()  # LTW!N29

$ws.Cells.Item(40, 8).Value = 33998.8  # LTW!H40 (was 19531.666)
$ws.Cells.Item(40, 9).Value = 41666.332  # LTW!I40 (was 18684.285)
$ws.Cells.Item(40, 11).Value = 41666.332  # LTW!K40 (was 18684.285)
$ws.Cells.Item(40, 13).Value = -41530.332  # LTW!M40 (was -18548.285)

$ws.Cells.Item(68, 8).Value = 2266.6667  # LTW!H68 (was 0)
$ws.Cells.Item(68, 9).Value = 1900  # LTW!I68 (was 0)
$ws.Cells.Item(68, 10).Value = 3000  # LTW!J68 (was 0)
$ws.Cells.Item(68, 11).Value = 1900  # LTW!K68 (was 0)
$ws.Cells.Item(68, 12).Value = 3000  # LTW!L68 (was 0)
$ws.Cells.Item(68, 13).Value = -1151  # LTW!M68 (was ADD)
$ws.Cells.Item(68, 14).Value = -4498  # LTW!N68 (was ADD)

$ws.Cells.Item(71, 8).Value = 2266.6667  # LTW!H71 (was 0)
$ws.Cells.Item(71, 9).Value = 1900  # LTW!I71 (was 0)
$ws.Cells.Item(71, 10).Value = 3000  # LTW!J71 (was 0)
$ws.Cells.Item(71, 11).Value = 9500  # LTW!K71 (was 0)
$ws.Cells.Item(71, 12).Value = 15000  # LTW!L71 (was 0)
$ws.Cells.Item(71, 13).Value = -5756  # LTW!M71 (was ADD)
$ws.Cells.Item(71, 14).Value = -22488  # LTW!N71 (was ADD)

$ws.Cells.Item(93, 8).Value = 28797.691  # LTW!H93 (was 31080.834)
$ws.Cells.Item(93, 9).Value = 3103  # LTW!I93 (was 3292.2222)
$ws.Cells.Item(93, 11).Value = 3103  # LTW!K93 (was 3292.2222)
$ws.Cells.Item(93, 13).Value = -1855  # LTW!M93 (was -2044.2222)

$ws.Cells.Item(105, 8).Value = 33871.668  # LTW!H105 (was 33538)
$ws.Cells.Item(105, 10).Value = 33871.668  # LTW!J105 (was 33538)
$ws.Cells.Item(105, 12).Value = 33871.668  # LTW!L105 (was 33538)
$ws.Cells.Item(105, 14).Value = -40859.668  # LTW!N105 (was -40526)

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 26599.4  # WVR!H62 (was 22866.166)
$ws.Cells.Item(62, 9).Value = 7665.6665  # WVR!I62 (was 6799.25)
$ws.Cells.Item(62, 11).Value = 7665.6665  # WVR!K62 (was 6799.25)
$ws.Cells.Item(62, 13).Value = -7041.6665  # WVR!M62 (was -6175.25)

$ws.Cells.Item(65, 8).Value = 26599.4  # WVR!H65 (was 22866.166)
$ws.Cells.Item(65, 9).Value = 7665.6665  # WVR!I65 (was 6799.25)
$ws.Cells.Item(65, 11).Value = 38328.3325  # WVR!K65 (was 33996.25)
$ws.Cells.Item(65, 13).Value = -35208.3325  # WVR!M65 (was -30876.25)

$ws.Cells.Item(113, 8).Value = 198.32259  # WVR!H113 (was 205.3)
$ws.Cells.Item(113, 9).Value = 177.13637  # WVR!I113 (was 198.04546)
$ws.Cells.Item(113, 10).Value = 250.11111  # WVR!J113 (was 225.25)
$ws.Cells.Item(113, 11).Value = 531.4091100000001  # WVR!K113 (was 594.1363799999999)
$ws.Cells.Item(113, 12).Value = 750.3333299999999  # WVR!L113 (was 675.75)
$ws.Cells.Item(113, 13).Value = 1638.59089  # WVR!M113 (was 1575.86362)
$ws.Cells.Item(113, 14).Value = -5090.333329999999  # WVR!N113 (was -5015.75)

$ws.Cells.Item(122, 8).Value = 2335.2  # WVR!H122 (was 2658.6667)
$ws.Cells.Item(122, 9).Value = 2335.2  # WVR!I122 (was 2658.6667)
$ws.Cells.Item(122, 11).Value = 7005.599999999999  # WVR!K122 (was 7976.000100000001)
$ws.Cells.Item(122, 13).Value = -4555.599999999999  # WVR!M122 (was -5526.000100000001)

$ws.Cells.Item(126, 8).Value = 1240.5  # WVR!H126 (was 1253.3125)
$ws.Cells.Item(126, 9).Value = 1179.0834  # WVR!I126 (was 1196.1666)
$ws.Cells.Item(126, 11).Value = 3537.2502  # WVR!K126 (was 3588.4998)
$ws.Cells.Item(126, 13).Value = -1067.2502  # WVR!M126 (was -1118.4998)

$ws.Cells.Item(132, 8).Value = 2022  # WVR!H132 (was 2055.125)
$ws.Cells.Item(132, 9).Value = 2101.75  # WVR!I132 (was 2147.0645)
$ws.Cells.Item(132, 11).Value = 6305.25  # WVR!K132 (was 6441.193499999999)
$ws.Cells.Item(132, 13).Value = -3775.25  # WVR!M132 (was -3911.193499999999)
